$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 212.75676
$ws.Range("I15").Value = 212.75676
$ws.Range("K15").Value = 638.2702800000001
$ws.Range("M15").Value = -469.2702800000001
$ws.Range("H32").Value = 5124.6875
$ws.Range("I32").Value = 3199.6
$ws.Range("J32").Value = 5999.727
$ws.Range("K32").Value = 3199.6
$ws.Range("L32").Value = 5999.727
$ws.Range("M32").Value = -2873.6
$ws.Range("N32").Value = -6651.727
$ws.Range("H33").Value = 156.66667
$ws.Range("I33").Value = 171
$ws.Range("K33").Value = 171
$ws.Range("M33").Value = 58
$ws.Range("H39").Value = 225.75
$ws.Range("I39").Value = 225.75
$ws.Range("K39").Value = 677.25
$ws.Range("M39").Value = -381.25
$ws.Range("H64").Value = 2172.25
$ws.Range("I64").Value = 1829.6666
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 1829.6666
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -1581.6666
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 2172.25
$ws.Range("I67").Value = 1829.6666
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 1829.6666
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -971.6666
$ws.Range("N67").Value = -4916
$ws.Range("H69").Value = 20506.5
$ws.Range("I69").Value = 20506.5
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 61519.5
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -60645.5
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 20506.5
$ws.Range("I72").Value = 20506.5
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 184558.5
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -180190.5
$ws.Range("N72").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H96").Value = 3012.652
$ws.Range("I96").Value = 2721.375
$ws.Range("J96").Value = 3678.4285
$ws.Range("K96").Value = 8164.125
$ws.Range("L96").Value = 11035.2855
$ws.Range("M96").Value = -6791.125
$ws.Range("N96").Value = -13781.2855
$ws.Range("H97").Value = 3345.5
$ws.Range("J97").Value = 3345.5
$ws.Range("L97").Value = 10036.5
$ws.Range("N97").Value = -11028.5
$ws.Range("H100").Value = 882.7143
$ws.Range("I100").Value = 882.7143
$ws.Range("K100").Value = 882.7143
$ws.Range("M100").Value = -341.7143
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H111").Value = 941.6667
$ws.Range("H115").Value = 690.5
$ws.Range("I115").Value = 1192.5714
$ws.Range("K115").Value = 3577.7142
$ws.Range("M115").Value = -2010.7142
$ws.Range("H137").Value = 5337.7036
$ws.Range("I137").Value = 6370.1
$ws.Range("K137").Value = 19110.3
$ws.Range("M137").Value = -16560.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8131.6665
$ws.Range("I31").Value = 6197.5
$ws.Range("J31").Value = 12000
$ws.Range("K31").Value = 6197.5
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = -5903.5
$ws.Range("N31").Value = -12588
$ws.Range("I34").Value = 15000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -14729
$ws.Range("N34").ClearContents()
$ws.Range("H45").Value = 10498.75
$ws.Range("I45").Value = 3999
$ws.Range("K45").Value = 3999
$ws.Range("M45").Value = -3622
$ws.Range("H61").Value = 9078
$ws.Range("I61").Value = 3924.6667
$ws.Range("K61").Value = 3924.6667
$ws.Range("M61").Value = -3712.6667
$ws.Range("H136").Value = 9078
$ws.Range("I136").Value = 3924.6667
$ws.Range("K136").Value = 11774.0001
$ws.Range("M136").Value = -9224.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 538.1
$ws.Range("J7").Value = 790.6
$ws.Range("L7").Value = 790.6
$ws.Range("N7").Value = -1016.6
$ws.Range("H22").Value = 389.75
$ws.Range("I22").Value = 343.375
$ws.Range("K22").Value = 343.375
$ws.Range("M22").Value = -170.375
$ws.Range("H23").Value = 750
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H94").Value = 711.4375
$ws.Range("I94").Value = 756.86664
$ws.Range("J94").Value = 30
$ws.Range("K94").Value = 756.86664
$ws.Range("L94").Value = 30
$ws.Range("M94").Value = -305.86664
$ws.Range("N94").Value = -932
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 1352.9
$ws.Range("I134").Value = 1352.9
$ws.Range("K134").Value = 4058.7
$ws.Range("M134").Value = -1523.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 389.51428
$ws.Range("I7").Value = 373.66666
$ws.Range("K7").Value = 373.66666
$ws.Range("M7").Value = -260.66666
$ws.Range("H22").Value = 33752
$ws.Range("I22").Value = 35686.65
$ws.Range("K22").Value = 35686.65
$ws.Range("M22").Value = -35336.65
$ws.Range("H62").Value = 30684.25
$ws.Range("J62").Value = 87099.39999999999
$ws.Range("L62").Value = 87099.39999999999
$ws.Range("N62").Value = -88347.39999999999
$ws.Range("H65").Value = 30684.25
$ws.Range("J65").Value = 87099.39999999999
$ws.Range("L65").Value = 435497
$ws.Range("N65").Value = -441737
$ws.Range("H86").Value = 28406.2
$ws.Range("I86").Value = 8093.2
$ws.Range("J86").Value = 48719.2
$ws.Range("K86").Value = 8093.2
$ws.Range("L86").Value = 48719.2
$ws.Range("M86").Value = -6970.2
$ws.Range("N86").Value = -50965.2
$ws.Range("H89").Value = 28406.2
$ws.Range("I89").Value = 8093.2
$ws.Range("J89").Value = 48719.2
$ws.Range("K89").Value = 40466
$ws.Range("L89").Value = 243596
$ws.Range("M89").Value = -34850
$ws.Range("N89").Value = -254828
$ws.Range("H92").Value = 78399.8
$ws.Range("J92").Value = 78399.8
$ws.Range("L92").Value = 78399.8
$ws.Range("N92").Value = -83391.8
$ws.Range("H132").Value = 3879.4
$ws.Range("I132").Value = 3754.889
$ws.Range("K132").Value = 11264.667
$ws.Range("M132").Value = -8734.667000000001
$ws.Range("H134").Value = 5097.7144
$ws.Range("I134").Value = 3424.75
$ws.Range("K134").Value = 10274.25
$ws.Range("M134").Value = -7739.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 3775.8333
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 3775.8333
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 11327.4999
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -12139.4999
$ws.Range("H85").Value = 3775.8333
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3775.8333
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 11327.4999
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -14135.4999
$ws.Range("H131").Value = 1546.7805
$ws.Range("J131").Value = 1550.4625
$ws.Range("L131").Value = 4651.387500000001
$ws.Range("N131").Value = -14731.3875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1207.7778
$ws.Range("I23").Value = 3700
$ws.Range("J23").Value = 495.7143
$ws.Range("K23").Value = 3700
$ws.Range("L23").Value = 495.7143
$ws.Range("M23").Value = -3477
$ws.Range("N23").Value = -941.7143
$ws.Range("H39").Value = 49997.5
$ws.Range("J39").Value = 49997.5
$ws.Range("L39").Value = 49997.5
$ws.Range("N39").Value = -51061.5
$ws.Range("H70").Value = 13398.148
$ws.Range("I70").Value = 13970.4
$ws.Range("K70").Value = 13970.4
$ws.Range("M70").Value = -13700.4
$ws.Range("H73").Value = 13398.148
$ws.Range("I73").Value = 13970.4
$ws.Range("K73").Value = 13970.4
$ws.Range("M73").Value = -13034.4
$ws.Range("H132").Value = 2559.6667
$ws.Range("I132").Value = 2559.6667
$ws.Range("K132").Value = 7679.000100000001
$ws.Range("M132").Value = -5149.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 25413.285
$ws.Range("I16").Value = 24648.834
$ws.Range("K16").Value = 24648.834
$ws.Range("M16").Value = -24478.834
$ws.Range("H35").Value = 2538.4614
$ws.Range("I35").Value = 1181.909
$ws.Range("J35").Value = 9999.5
$ws.Range("K35").Value = 1181.909
$ws.Range("L35").Value = 9999.5
$ws.Range("M35").Value = -845.9090000000001
$ws.Range("N35").Value = -10671.5
$ws.Range("H55").Value = 2482.4
$ws.Range("J55").Value = 5494.778
$ws.Range("L55").Value = 5494.778
$ws.Range("N55").Value = -5840.778
$ws.Range("H61").Value = 28581.666
$ws.Range("I61").Value = 26298.4
$ws.Range("K61").Value = 26298.4
$ws.Range("M61").Value = -26096.4
$ws.Range("H68").Value = 22666.166
$ws.Range("I68").Value = 9000.25
$ws.Range("K68").Value = 9000.25
$ws.Range("M68").Value = -8251.25
$ws.Range("H71").Value = 22666.166
$ws.Range("I71").Value = 9000.25
$ws.Range("K71").Value = 45001.25
$ws.Range("M71").Value = -41257.25
$ws.Range("H93").Value = 24998
$ws.Range("I93").Value = 24998
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 24998
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -23750
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 87998
$ws.Range("J94").Value = 87998
$ws.Range("L94").Value = 87998
$ws.Range("N94").Value = -89350
$ws.Range("H101").Value = 20724.5
$ws.Range("J101").Value = 20724.5
$ws.Range("L101").Value = 20724.5
$ws.Range("N101").Value = -27214.5
$ws.Range("H113").Value = 28581.666
$ws.Range("I113").Value = 26298.4
$ws.Range("K113").Value = 26298.4
$ws.Range("M113").Value = -24128.4
$ws.Range("H122").Value = 16130.667
$ws.Range("I122").Value = 16130.667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 48392.001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -45942.001
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 42497.25
$ws.Range("I132").Value = 39997
$ws.Range("K132").Value = 119991
$ws.Range("M132").Value = -117461
$ws.Range("H136").Value = 10817.637
$ws.Range("J136").Value = 18732.666
$ws.Range("L136").Value = 56197.99800000001
$ws.Range("N136").Value = -61297.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 9750
$ws.Range("J53").Value = 9500
$ws.Range("L53").Value = 9500
$ws.Range("N53").Value = -10714
$ws.Range("H58").Value = 20407.637
$ws.Range("I58").Value = 12914.667
$ws.Range("J58").Value = 29399.2
$ws.Range("K58").Value = 12914.667
$ws.Range("L58").Value = 29399.2
$ws.Range("M58").Value = -12606.667
$ws.Range("N58").Value = -30015.2
$ws.Range("H63").Value = 13848.8
$ws.Range("J63").Value = 16622
$ws.Range("L63").Value = 16622
$ws.Range("N63").Value = -17870
$ws.Range("H66").Value = 13848.8
$ws.Range("J66").Value = 16622
$ws.Range("L66").Value = 49866
$ws.Range("N66").Value = -56106
$ws.Range("H82").Value = 30150
$ws.Range("I82").Value = 28000
$ws.Range("K82").Value = 28000
$ws.Range("M82").Value = -27617
$ws.Range("H85").Value = 30150
$ws.Range("I85").Value = 28000
$ws.Range("K85").Value = 28000
$ws.Range("M85").Value = -26674
$ws.Range("H100").Value = 1156.4
$ws.Range("I100").Value = 933.2222
$ws.Range("J100").Value = 1491.1666
$ws.Range("K100").Value = 1866.4444
$ws.Range("L100").Value = 2982.3332
$ws.Range("M100").Value = -1325.4444
$ws.Range("N100").Value = -4064.3332
$ws.Range("H103").Value = 45344.832
$ws.Range("J103").Value = 45344.832
$ws.Range("L103").Value = 45344.832
$ws.Range("N103").Value = -47688.832
$ws.Range("H107").Value = 1244.826
$ws.Range("I107").Value = 1085.4615
$ws.Range("K107").Value = 3256.3845
$ws.Range("M107").Value = -1336.3845
$ws.Range("H122").Value = 2150.7273
$ws.Range("I122").Value = 1930.5714
$ws.Range("K122").Value = 5791.7142
$ws.Range("M122").Value = -3341.7142
$ws.Range("H132").Value = 2998.5334
$ws.Range("I132").Value = 2844.4614
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8533.3842
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6003.3842
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 11202.667
$ws.Range("I136").Value = 4140.9287
$ws.Range("K136").Value = 12422.7861
$ws.Range("M136").Value = -9872.786100000001
